{"js": "// Replace the \"Phone: 612-840-9988 \" paragraph text with \"Email: rvksrinivas@yahoo.com\"\nconst body = context.document.body;\nconst results = body.search(\"Phone: 612-840-9988 \", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  const range = results.items[i];\n  const para = range.paragraphs.getFirstOrNullObject();\n  para.load(\"text\");\n  await context.sync();\n\n  // Clear the paragraph's contents and insert the new text fresh, so the\n  // resulting run has no leftover/preserved whitespace formatting.\n  para.clear();\n  await context.sync();\n\n  para.insertText(\"Email: rvksrinivas@yahoo.com\", Word.InsertLocation.start);\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"Phone: 612-840-9988 \"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"Email: rvksrinivas@yahoo.com\"\n$find.Execute(\n    $find.Text,       # FindText\n    $false,           # MatchCase\n    $false,           # MatchWholeWord\n    $false,           # MatchWildcards\n    $false,           # MatchSoundsLike\n    $false,           # MatchAllWordForms\n    $true,            # Forward\n    1,                # Wrap (wdFindContinue)\n    $false,           # Format\n    $find.Replacement.Text,  # ReplaceWith\n    2                 # Replace (wdReplaceAll)\n)\n"}
